$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(194, "Batman: Arkham Asylum Batmobile", "Mainlines")
    ,@(195, "TV Series Batmobile", "Mainlines")
    ,@(196, "The Bat", "Mainlines")
    ,@(197, "Ford GT-40", "Silver Series National Icons")
    ,@(198, "Porsche 914 Safari", "HW Speed Graphics")
    ,@(199, "16 Lamborghini Centenario Roadster", "Mainlines")
    ,@(200, "2018 Bentley Continental GT3", "Premiums Fast & Furious")
    ,@(201, "Mazda RX7 FC Pandem", "Premiums Car Culture")
    ,@(202, "04 Mazda Mazdaspeed Miata", "Premiums Car Culture")
    ,@(203, "Audi RS E-Tron GT", "Mainlines")
    ,@(204, "Maserati MC20", "Premiums Boulevard")
    ,@(205, "1989 Batmobile", "Mainlines")
    ,@(206, "Justice League Batmobile", "Mainlines")
    ,@(207, "Custom Cadillac Fleetwood", "53rd Anniversary Series")
    ,@(208, "Corvette Grand Sport Roadster", "Mainlines")
    ,@(209, "Porsche 911 Carrera T", "Mainlines")
    ,@(210, "37 Bugatti", "Mainlines")
    ,@(211, "Bentley Continental Supersports", "Exotics")
    ,@(212, "BMW M3 Wagon", "Mainlines")
    ,@(213, "Volvo 240 Drift Wagon", "Mainlines")
    ,@(214, "2020 Dodge Charger `"Hellcat`"", "Mainlines")
)

$startRow = 195
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
